$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.124.67'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.546.92'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.89'
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.52'
$ws.Range("E6").Value = '  -4.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.547.18'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("E10").Value = '  -2.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.08'
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.387'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.151.72'
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("E14").Value = '  -3.09%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.551.92'
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.75'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.286.39'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.93'
$ws.Range("E19").Value = '  -3.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.34'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.39'
$ws.Range("E22").Value = '  -1.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.576'
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.694.76'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.91'
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.71'
$ws.Range("E28").Value = '  -0.54%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.49'
$ws.Range("E30").Value = '  +2.68%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.53'
$ws.Range("E31").Value = '  +23.67%  '
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.555.59'
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.05'
$ws.Range("E34").Value = '  +0.59%  '
$ws.Range("E36").Value = '  +0.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '169.95'
$ws.Range("E37").Value = '  +0.58%  '
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.99'
$ws.Range("E40").Value = '  +1.34%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.20'
$ws.Range("E43").Value = '  -1.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.04'
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.23'
$ws.Range("E45").Value = '  +3.66%  '
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.456.94'
$ws.Range("E49").Value = '  +3.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.89'
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0264'
$ws.Range("E51").Value = '  +1.15%  '
